# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted at row 118 (right before the
# existing row that used to be there), pushing all subsequent rows
# (old 118..165) down by one to (119..166) and extending the used range
# from A1:R165 to A1:R166.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 118..165 down to 119..166, carrying formatting along
# (mirrors Excel's own "insert row" behaviour).
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new observation.
$ws.Cells.Item(118, 1).Value = 8
$ws.Cells.Item(118, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44489
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 100112012
$ws.Cells.Item(118, 7).Value = "Espinaca"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 3000
$ws.Cells.Item(118, 11).Value = 450
$ws.Cells.Item(118, 12).Value = 500
$ws.Cells.Item(118, 13).Value = 475
$ws.Cells.Item(118, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(118, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(118, 16).Value = 950
$ws.Cells.Item(118, 17).Value = 0.5
$ws.Cells.Item(118, 18).Value = "Hortaliza"
